$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vanilla non-needy modules added to the reference sheet (Steel Crate Games, 2015-10-08)
$rows = @(
    @("Wires", "Wires", 1, "modules/Wires.pdf", "Steel Crate Games", "2015-10-08", 0),
    @("Button", "BigButton", 1, "modules/The Button.pdf", "Steel Crate Games", "2015-10-08", 0),
    @("Keypad", "Keypad", 2, "modules/Keypad.pdf", "Steel Crate Games", "2015-10-08", 0),
    @("Simon Says", "Simon", 1, "modules/Simon Says.pdf", "Steel Crate Games", "2015-10-08", 0),
    @("Who's On First", "WhosOnFirst", 2, "modules/Who's On First.pdf", "Steel Crate Games", "2015-10-08", 0),
    @("Memory", "Memory", 2, "modules/Memory.pdf", "Steel Crate Games", "2015-10-08", 0),
    @("Morse Code", "Morse", 3, "modules/Morse Code.pdf", "Steel Crate Games", "2015-10-08", 0),
    @("Complicated Wires", "Venn", 3, "modules/Complicated Wires.pdf", "Steel Crate Games", "2015-10-08", 0),
    @("Wire Sequence", "WireSequence", 3, "modules/Wire Sequence.pdf", "Steel Crate Games", "2015-10-08", 0),
    @("Maze", "Maze", 2, "modules/Maze.pdf", "Steel Crate Games", "2015-10-08", 0),
    @("Password", "Password", 2, "modules/Password.pdf", "Steel Crate Games", "2015-10-08", 0)
)

$startRow = 99
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

# Update the window/view state to reflect scrolling down to the newly added rows
[void]$ws.Range("A110").Select()
